# Apply Diebold-Mariano summary correction
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MCPS): Comparaciones_Significativas 6/10 -> 5/10, ECRPS_Mejor 153.6 -> 128
$ws.Range("B2").Value = "5/10"
$ws.Range("C2").Value = 128

# Row 3: model label swaps from AV-MCPS -> Sieve Bootstrap;
# Mejor_N_Calib 200 -> 100, ECRPS_Mejor 0.6217869373211101 -> 0.5295197719051405
$ws.Range("A3").Value = "Sieve Bootstrap"
$ws.Range("D3").Value = 100
$ws.Range("E3").Value = 0.5295197719051405

# Row 4: model label swaps from Sieve Bootstrap -> AV-MCPS;
# Mejor_N_Calib 100 -> 200, ECRPS_Mejor 0.5295197719051405 -> 0.6217869373211101
$ws.Range("A4").Value = "AV-MCPS"
$ws.Range("D4").Value = 200
$ws.Range("E4").Value = 0.6217869373211101

$wb.Save()
